$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 16.44017404108616
$ws.Range("C2").Value = 0.0000000000000001110223024625157
$ws.Range("D2").Value = 0.01533053832641091
$ws.Range("E2").Value = 0.9287579776373168
$ws.Range("F2").Value = 0.8625913810249587

$ws.Range("B3").Value = 16.46040920501808
$ws.Range("C3").Value = 0.0000000000000001110223024625157
$ws.Range("D3").Value = 0.01530198511979786
$ws.Range("E3").Value = 0.9270281611191774
$ws.Range("F3").Value = 0.8593812115080036

$ws.Range("B4").Value = 17.32862426996585
$ws.Range("C4").Value = 0.0000000000000001110223024625157
$ws.Range("D4").Value = 0.01299564749400441
$ws.Range("E4").Value = 0.7873051179047998
$ws.Range("F4").Value = 0.6198493486790908

$ws.Range("B5").Value = 17.66181136376285
$ws.Range("C5").Value = 0.0000000000000001110223024625157
$ws.Range("D5").Value = 0.0134736126112655
$ws.Range("E5").Value = 0.8162613036718603
$ws.Range("F5").Value = 0.6662825158720849

$ws.Range("B6").Value = 17.32696186820391
$ws.Range("C6").Value = 0.0000000000000001110223024625157
$ws.Range("D6").Value = 0.01297666684029057
$ws.Range("E6").Value = 0.7861552278499201
$ws.Range("F6").Value = 0.6180400422757597
